$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update CoarseSubstrate_score for row 2 (J2): 5 -> 3
$ws.Range("J2").Value = 3

# Add PoolQuantity&Quality_score values in column O for rows 2-5
$ws.Range("O2").Value = 1
$ws.Range("O3").Value = 5
$ws.Range("O4").Value = 5
$ws.Range("O5").Value = 5
